$wb = $excel.ActiveWorkbook

# sheet1 (展览)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1804
$ws.Range("G3").Value = 80
$ws.Range("F4").Value = 476
$ws.Range("G5").Value = "不可售"
$ws.Range("F7").Value = 640
$ws.Range("G7").Value = 68
$ws.Range("F8").Value = 354
$ws.Range("F9").Value = 1766
$ws.Range("F10").Value = 383
$ws.Range("F11").Value = 1440
$ws.Range("F12").Value = 826
$ws.Range("F13").Value = 351
$ws.Range("F14").Value = 696
$ws.Range("F15").Value = 12931
$ws.Range("F16").Value = 12886
$ws.Range("F17").Value = 966
$ws.Range("F20").Value = 529
$ws.Range("F21").Value = 57
$ws.Range("F22").Value = 597
$ws.Range("F23").Value = 2023
$ws.Range("F25").Value = 16
$ws.Range("F26").Value = 9
$ws.Range("F28").Value = 110
$ws.Range("F29").Value = 259
$ws.Range("F30").Value = 693

# sheet2 (演出)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 19
$ws.Range("F7").Value = 13

# sheet3 (本地生活)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 87
$ws.Range("F3").Value = 178

# sheet4 (全部类型)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 87
$ws.Range("F3").Value = 178
$ws.Range("F5").Value = 1804
$ws.Range("G5").Value = 80
$ws.Range("F6").Value = 476
$ws.Range("G8").Value = "不可售"
$ws.Range("F11").Value = 640
$ws.Range("G11").Value = 68
$ws.Range("F13").Value = 354
$ws.Range("F14").Value = 1766
$ws.Range("F15").Value = 383
$ws.Range("F16").Value = 1440
$ws.Range("F17").Value = 826
$ws.Range("F18").Value = 351
$ws.Range("F20").Value = 696
$ws.Range("F21").Value = 12931
$ws.Range("F22").Value = 12886
$ws.Range("F23").Value = 966
$ws.Range("F26").Value = 529
$ws.Range("F27").Value = 57
$ws.Range("F28").Value = 597
$ws.Range("F29").Value = 19
$ws.Range("F30").Value = 13
$ws.Range("F31").Value = 2023
$ws.Range("F33").Value = 16
$ws.Range("F34").Value = 9
$ws.Range("F38").Value = 110
$ws.Range("F39").Value = 259
$ws.Range("F40").Value = 693
